$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 248
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:27:03'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = 'Atasco tuerca'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:27:05'
$row[0,8] = '0:00:02'
$ws.Range("A248:I248").Value = $row

# Row 249
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:32:36'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = 'No coloca bien el sealling'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:33:41'
$row[0,8] = '0:01:05'
$ws.Range("A249:I249").Value = $row

# Row 250
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:34:21'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = 'La cámara no detecta Busbar'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:34:33'
$row[0,8] = '0:00:12'
$ws.Range("A250:I250").Value = $row

# Row 251
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:35:20'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = 'No coloca bien el sealling'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:35:21'
$row[0,8] = '0:00:01'
$ws.Range("A251:I251").Value = $row

# Row 252
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:40:07'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = 'Detección de sealling mal puesto'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:40:09'
$row[0,8] = '0:00:02'
$ws.Range("A252:I252").Value = $row

# Row 253
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:40:10'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = 'No lee QR'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:40:11'
$row[0,8] = '0:00:01'
$ws.Range("A253:I253").Value = $row

# Row 254
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:40:50'
$row[0,2] = '-'
$row[0,3] = 'Cámara no detecta skeleton'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:40:52'
$row[0,8] = '0:00:02'
$ws.Range("A254:I254").Value = $row

# Row 255
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:41:05'
$row[0,2] = '-'
$row[0,3] = 'Etiquetadora'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:41:07'
$row[0,8] = '0:00:02'
$ws.Range("A255:I255").Value = $row

# Row 256
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:42:11'
$row[0,2] = '-'
$row[0,3] = 'Robot no coge busbar'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:42:13'
$row[0,8] = '0:00:02'
$ws.Range("A256:I256").Value = $row

# Row 257
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:45:30'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = 'Traza'
$row[0,6] = '-'
$row[0,7] = '09:45:32'
$row[0,8] = '0:00:02'
$ws.Range("A257:I257").Value = $row

# Row 258
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:45:34'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = 'No coloca bien la pcb'
$row[0,6] = '-'
$row[0,7] = '09:45:35'
$row[0,8] = '0:00:01'
$ws.Range("A258:I258").Value = $row

# Row 259
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:45:39'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = 'Fallo atornillador'
$row[0,6] = '-'
$row[0,7] = '09:45:41'
$row[0,8] = '0:00:02'
$ws.Range("A259:I259").Value = $row

# Row 260
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:46:53'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = 'NOK Soldad. Plástico+Metal'
$row[0,6] = '-'
$row[0,7] = '09:46:56'
$row[0,8] = '0:00:03'
$ws.Range("A260:I260").Value = $row

# Row 261
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:47:37'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = 'Fallo cámara ferrite'
$row[0,6] = '-'
$row[0,7] = '09:47:39'
$row[0,8] = '0:00:02'
$ws.Range("A261:I261").Value = $row

# Row 262
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:49:35'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = 'Robot no coge PCB'
$row[0,6] = '-'
$row[0,7] = '09:49:37'
$row[0,8] = '0:00:02'
$ws.Range("A262:I262").Value = $row

# Row 263
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:50:22'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = 'NOK Soldad. Plástico+Metal'
$row[0,6] = '-'
$row[0,7] = '09:50:25'
$row[0,8] = '0:00:03'
$ws.Range("A263:I263").Value = $row

# Row 264
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:53:36'
$row[0,2] = '-'
$row[0,3] = 'Tornillo atascado en tolva'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:53:43'
$row[0,8] = '0:00:07'
$ws.Range("A264:I264").Value = $row

# Row 265
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:53:44'
$row[0,2] = '-'
$row[0,3] = 'Detección de sealling mal puesto'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:53:46'
$row[0,8] = '0:00:02'
$ws.Range("A265:I265").Value = $row

# Row 266
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:54:11'
$row[0,2] = '-'
$row[0,3] = 'Cámara no detecta Pcb'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:54:12'
$row[0,8] = '0:00:01'
$ws.Range("A266:I266").Value = $row

# Row 267
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '09:54:40'
$row[0,2] = '-'
$row[0,3] = 'Etiquetadora'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '09:54:41'
$row[0,8] = '0:00:01'
$ws.Range("A267:I267").Value = $row

# Row 268
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '10:40:39'
$row[0,2] = 'No atornilla clips'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = '-'
$row[0,7] = '10:40:40'
$row[0,8] = '0:00:01'
$ws.Range("A268:I268").Value = $row

# Row 269
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '11:02:33'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = 'Fallo dispensación glue'
$row[0,7] = '11:02:36'
$row[0,8] = '0:00:03'
$ws.Range("A269:I269").Value = $row

# Row 270
$row = New-Object 'object[,]' 1,9
$row[0,0] = '''2024-05-24'
$row[0,1] = '11:02:35'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = '-'
$row[0,6] = 'Error en sensor de salida'
$row[0,7] = '11:02:36'
$row[0,8] = '0:00:01'
$ws.Range("A270:I270").Value = $row
